# BBSalary: output *all* wage details (for every selected employee),
# not just the first one. This adds a second employee column (C) that
# mirrors column B, and promotes the employee name up into the header
# row (row 1) instead of a separate "Name" label row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old layout had a standalone "Name" row (row 2: A2="Name",
# B2=<employee name>). Remove it entirely - row 3 ("Gross Salary")
# and everything below shifts up to become the new row 2, etc.
$ws.Rows.Item(2).Delete()

# Put the first employee's name directly into the header cell (B1),
# keeping the bold/boxed header style that was already there.
$ws.Range("B1").Value = "Janet Apostol"

# Add a second employee column. Clone B1's header style into C1 so it
# matches (bold, bordered, centered), then set its name.
$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("C1").Value = "Rizalyn Repalda"

# Mirror the first employee's wage detail values (Gross Salary,
# PhilHealth/SSS/Pag-Ibig Deductions, Net Salary) into column C for
# the second employee, carrying over the (unstyled) number formatting.
$ws.Range("B2").Copy($ws.Range("C2"))
$ws.Range("B3").Copy($ws.Range("C3"))
$ws.Range("B4").Copy($ws.Range("C4"))
$ws.Range("B5").Copy($ws.Range("C5"))
$ws.Range("B6").Copy($ws.Range("C6"))
